$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049656380268245
$ws.Range("D2").Value = 1.058211182012943
$ws.Range("E2").Value = 1.0568314220712
$ws.Range("F2").Value = 1.067561923627415
$ws.Range("I2").Value = 1.04916475645157
$ws.Range("J2").Value = 1.054693671362677
$ws.Range("K2").Value = 1.060944036469356
$ws.Range("L2").Value = 1.059568054282428
$ws.Range("M2").Value = 1.070269454366599
$ws.Range("N2").Value = 1.022006308103492

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050716512279256
$ws.Range("D3").Value = 1.058933587860115
$ws.Range("E3").Value = 1.057746352548941
$ws.Range("F3").Value = 1.068493663779592
$ws.Range("I3").Value = 1.049459089842716
$ws.Range("J3").Value = 1.055402629931996
$ws.Range("K3").Value = 1.061480663106086
$ws.Range("L3").Value = 1.060296443913907
$ws.Range("M3").Value = 1.071016715151169
$ws.Range("N3").Value = 1.022246933887035

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051402771337858
$ws.Range("D4").Value = 1.059401156220499
$ws.Range("E4").Value = 1.058338952975786
$ws.Range("F4").Value = 1.069097149255698
$ws.Range("I4").Value = 1.049648310768965
$ws.Range("J4").Value = 1.055861062871362
$ws.Range("K4").Value = 1.061827347479925
$ws.Range("L4").Value = 1.060767707128274
$ws.Range("M4").Value = 1.071500198686111
$ws.Range("N4").Value = 1.022402410313867

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051691342159222
$ws.Range("D5").Value = 1.059597750409831
$ws.Range("E5").Value = 1.058588220051641
$ws.Range("F5").Value = 1.069350994188324
$ws.Range("I5").Value = 1.049727563790327
$ws.Range("J5").Value = 1.05605371321758
$ws.Range("K5").Value = 1.06197296142384
$ws.Range("L5").Value = 1.06096581273358
$ws.Range("M5").Value = 1.071703443814323
$ws.Range("N5").Value = 1.02246771858792

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051739798423794
$ws.Range("D6").Value = 1.059630761075369
$ws.Range("E6").Value = 1.058630081150895
$ws.Range("F6").Value = 1.069393624030031
$ws.Range("I6").Value = 1.049740853383979
$ws.Range("J6").Value = 1.056086055654344
$ws.Range("K6").Value = 1.061997402873659
$ws.Range("L6").Value = 1.06099907472965
$ws.Range("M6").Value = 1.071737568877515
$ws.Range("N6").Value = 1.022478680959466

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051406626970229
$ws.Range("D7").Value = 1.059403783011283
$ws.Range("E7").Value = 1.058342283153988
$ws.Range("F7").Value = 1.069100540596633
$ws.Range("I7").Value = 1.049649370913157
$ws.Range("J7").Value = 1.055863637369481
$ws.Range("K7").Value = 1.061829293700836
$ws.Range("L7").Value = 1.060770354278696
$ws.Range("M7").Value = 1.071502914503035
$ws.Range("N7").Value = 1.02240328317879

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050014598188683
$ws.Range("D8").Value = 1.058455296219865
$ws.Range("E8").Value = 1.057140506688586
$ws.Range("F8").Value = 1.067876687393519
$ws.Range("I8").Value = 1.049264482745624
$ws.Range("J8").Value = 1.054933331261091
$ws.Range("K8").Value = 1.061125505166156
$ws.Range("L8").Value = 1.059814227689291
$ws.Range("M8").Value = 1.07052200332359
$ws.Range("N8").Value = 1.022087675051393

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047563838853953
$ws.Range("D9").Value = 1.056784935700791
$ws.Range("E9").Value = 1.055027294845132
$ws.Range("F9").Value = 1.065724641869396
$ws.Range("I9").Value = 1.048576839136172
$ws.Range("J9").Value = 1.05329165509129
$ws.Range("K9").Value = 1.059881168486438
$ws.Range("L9").Value = 1.058129027364434
$ws.Range("M9").Value = 1.068793208420777
$ws.Range("N9").Value = 1.021529824828717

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04593146806375
$ws.Range("D10").Value = 1.055672095231603
$ws.Range("E10").Value = 1.053621539961588
$ws.Range("F10").Value = 1.064293052769908
$ws.Range("I10").Value = 1.048112098555884
$ws.Range("J10").Value = 1.05219564607019
$ws.Range("K10").Value = 1.059048849376609
$ws.Range("L10").Value = 1.057005338141598
$ws.Range("M10").Value = 1.067640518623316
$ws.Range("N10").Value = 1.021156790859933

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045224982579521
$ws.Range("D11").Value = 1.055190410213152
$ws.Range("E11").Value = 1.05301356596378
$ws.Range("F11").Value = 1.063673907209321
$ws.Range("I11").Value = 1.047909368289318
$ws.Range("J11").Value = 1.051720697631652
$ws.Range("K11").Value = 1.058687799831713
$ws.Range("L11").Value = 1.056518722499786
$ws.Range("M11").Value = 1.067141362226947
$ws.Range("N11").Value = 1.020994997049184

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044962613720388
$ws.Range("D12").Value = 1.055011519360587
$ws.Range("E12").Value = 1.052787847280828
$ws.Range("F12").Value = 1.063444041295728
$ws.Range("I12").Value = 1.047833841019937
$ws.Range("J12").Value = 1.051544225494213
$ws.Range("K12").Value = 1.05855359297479
$ws.Range("L12").Value = 1.056337964648175
$ws.Range("M12").Value = 1.06695594885209
$ws.Range("N12").Value = 1.020934859626611

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04501889041038
$ws.Range("D13").Value = 1.055049890753589
$ws.Range("E13").Value = 1.052836259718982
$ws.Range("F13").Value = 1.063493343218465
$ws.Range("I13").Value = 1.047850052017052
$ws.Range("J13").Value = 1.051582081856433
$ws.Range("K13").Value = 1.058582385179196
$ws.Range("L13").Value = 1.056376738135882
$ws.Range("M13").Value = 1.066995720841817
$ws.Range("N13").Value = 1.020947761114445

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045203294041432
$ws.Range("D14").Value = 1.055175622453966
$ws.Range("E14").Value = 1.052994905740535
$ws.Range("F14").Value = 1.063654904125129
$ws.Range("I14").Value = 1.047903129753904
$ws.Range("J14").Value = 1.051706111512576
$ws.Range("K14").Value = 1.058676708225625
$ws.Range("L14").Value = 1.056503781142081
$ws.Range("M14").Value = 1.067126035986794
$ws.Range("N14").Value = 1.020990026885778

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045316918001697
$ws.Range("D15").Value = 1.055253093699638
$ws.Range("E15").Value = 1.053092667377666
$ws.Range("F15").Value = 1.063754462039959
$ws.Range("I15").Value = 1.047935802997227
$ws.Range("J15").Value = 1.051782522980888
$ws.Range("K15").Value = 1.058734810936397
$ws.Range("L15").Value = 1.056582055600441
$ws.Range("M15").Value = 1.067206326876742
$ws.Range("N15").Value = 1.021016062928765

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045978362376241
$ws.Range("D16").Value = 1.055704067007582
$ws.Range("E16").Value = 1.053661904572023
$ws.Range("F16").Value = 1.064334159125995
$ws.Range("I16").Value = 1.048125521634784
$ws.Range("J16").Value = 1.052227159082349
$ws.Range("K16").Value = 1.059072797409687
$ws.Range("L16").Value = 1.057037632191224
$ws.Range("M16").Value = 1.067673645327129
$ws.Range("N16").Value = 1.021167522961557

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046393360439769
$ws.Range("D17").Value = 1.05598700016975
$ws.Range("E17").Value = 1.054019167185271
$ws.Range("F17").Value = 1.06469798703339
$ws.Range("I17").Value = 1.048244127217486
$ws.Range("J17").Value = 1.052505968931179
$ws.Range("K17").Value = 1.059284633867827
$ws.Range("L17").Value = 1.057323390101183
$ws.Range("M17").Value = 1.067966773179009
$ws.Range("N17").Value = 1.021262458333645

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046635454774833
$ws.Range("D18").Value = 1.056152047811806
$ws.Range("E18").Value = 1.054227622536019
$ws.Range("F18").Value = 1.064910273228791
$ws.Range("I18").Value = 1.048313163625748
$ws.Range("J18").Value = 1.052668558343751
$ws.Range("K18").Value = 1.059408131742073
$ws.Range("L18").Value = 1.057490062921349
$ws.Range("M18").Value = 1.068137746360536
$ws.Range("N18").Value = 1.021317806692301

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046718008283413
$ws.Range("D19").Value = 1.056208327732092
$ws.Range("E19").Value = 1.054298712295181
$ws.Range("F19").Value = 1.064982669442956
$ws.Range("I19").Value = 1.048336678784903
$ws.Range("J19").Value = 1.052723991073612
$ws.Range("K19").Value = 1.059450230662622
$ws.Range("L19").Value = 1.057546893193852
$ws.Range("M19").Value = 1.06819604321132
$ws.Range("N19").Value = 1.021336674663285

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.046348831663444
$ws.Range("D20").Value = 1.055956642299923
$ws.Range("E20").Value = 1.053980829030859
$ws.Range("F20").Value = 1.064658944337069
$ws.Range("I20").Value = 1.048231416876136
$ws.Range("J20").Value = 1.052476058958468
$ws.Range("K20").Value = 1.059261912306183
$ws.Range("L20").Value = 1.057292731488547
$ws.Range("M20").Value = 1.067935323674783
$ws.Range("N20").Value = 1.021252275337638

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045148990352321
$ws.Range("D21").Value = 1.055138596809465
$ws.Range("E21").Value = 1.05294818540123
$ws.Range("F21").Value = 1.0636073253577
$ws.Range("I21").Value = 1.047887505870417
$ws.Range("J21").Value = 1.051669589389579
$ws.Range("K21").Value = 1.058648935111249
$ws.Range("L21").Value = 1.056466370336214
$ws.Range("M21").Value = 1.067087661542553
$ws.Range("N21").Value = 1.020977581771983

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04439489931713
$ws.Range("D22").Value = 1.054624424022851
$ws.Range("E22").Value = 1.05229955772567
$ws.Range("F22").Value = 1.062946781004428
$ws.Range("I22").Value = 1.047669978249805
$ws.Range("J22").Value = 1.051162211109395
$ws.Range("K22").Value = 1.05826297124702
$ws.Range("L22").Value = 1.055946763242853
$ws.Range("M22").Value = 1.06655467721433
$ws.Range("N22").Value = 1.020804639717573

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.044794629123764
$ws.Range("D23").Value = 1.054896980773475
$ws.Range("E23").Value = 1.052643346980541
$ws.Range("F23").Value = 1.063296886147519
$ws.Range("I23").Value = 1.047785416604062
$ws.Range("J23").Value = 1.051431212044315
$ws.Range("K23").Value = 1.058467630929754
$ws.Range("L23").Value = 1.056222220518626
$ws.Range("M23").Value = 1.066837224487968
$ws.Range("N23").Value = 1.020896341426984

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046368952198683
$ws.Range("D24").Value = 1.055970359660516
$ws.Range("E24").Value = 1.053998152175959
$ws.Range("F24").Value = 1.06467658582997
$ws.Range("I24").Value = 1.048237160577891
$ws.Range("J24").Value = 1.052489574097799
$ws.Range("K24").Value = 1.059272179395758
$ws.Range("L24").Value = 1.057306584811141
$ws.Range("M24").Value = 1.067949534362328
$ws.Range("N24").Value = 1.021256876675045

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048197160577859
$ws.Range("D25").Value = 1.05721663942323
$ws.Range("E25").Value = 1.055573076044033
$ws.Range("F25").Value = 1.066280453481247
$ws.Range("I25").Value = 1.048755725624007
$ws.Range("J25").Value = 1.053716344372003
$ws.Range("K25").Value = 1.060203349211975
$ws.Range("L25").Value = 1.058564733814121
$ws.Range("M25").Value = 1.069240174782798
$ws.Range("N25").Value = 1.021674243410187
